# Auto-generated edit script: update Hyperion Profits market-price data
# Mirrors a scheduled runner pulling fresh Universalis prices into the
# per-job (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) "Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4176.189
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 4368.206
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 4368.206
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -5020.206
$ws.Range("H62").Value = 7899.15
$ws.Range("I62").Value = 997
$ws.Range("K62").Value = 997
$ws.Range("M62").Value = -373
$ws.Range("H65").Value = 7899.15
$ws.Range("I65").Value = 997
$ws.Range("K65").Value = 4985
$ws.Range("M65").Value = -1865
$ws.Range("H74").Value = 6242
$ws.Range("I74").Value = 3330.75
$ws.Range("K74").Value = 3330.75
$ws.Range("M74").Value = -2394.75
$ws.Range("H76").Value = 2530799.8
$ws.Range("I76").Value = 4278233.5
$ws.Range("J76").Value = 6729.222
$ws.Range("K76").Value = 4278233.5
$ws.Range("L76").Value = 6729.222
$ws.Range("M76").Value = -4277918.5
$ws.Range("N76").Value = -7359.222
$ws.Range("H77").Value = 6242
$ws.Range("I77").Value = 3330.75
$ws.Range("K77").Value = 16653.75
$ws.Range("M77").Value = -11973.75
$ws.Range("H79").Value = 2530799.8
$ws.Range("I79").Value = 4278233.5
$ws.Range("J79").Value = 6729.222
$ws.Range("K79").Value = 4278233.5
$ws.Range("L79").Value = 6729.222
$ws.Range("M79").Value = -4277141.5
$ws.Range("N79").Value = -8913.222
$ws.Range("H86").Value = 5243.875
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 5992.5713
$ws.Range("K86").Value = 3
$ws.Range("L86").Value = 5992.5713
$ws.Range("M86").Value = 1120
$ws.Range("N86").Value = -8238.5713
$ws.Range("H89").Value = 5243.875
$ws.Range("I89").Value = 3
$ws.Range("J89").Value = 5992.5713
$ws.Range("K89").Value = 15
$ws.Range("L89").Value = 29962.8565
$ws.Range("M89").Value = 5601
$ws.Range("N89").Value = -41194.85649999999
$ws.Range("H112").Value = 4524.9443
$ws.Range("J112").Value = 4524.9443
$ws.Range("L112").Value = 13574.8329
$ws.Range("N112").Value = -15790.8329
$ws.Range("H132").Value = 1991
$ws.Range("I132").Value = 2023.5476
$ws.Range("K132").Value = 6070.642800000001
$ws.Range("M132").Value = -3540.642800000001
$ws.Range("H135").Value = 1745.2258
$ws.Range("I135").Value = 814.0952
$ws.Range("K135").Value = 7326.8568
$ws.Range("M135").Value = -4791.8568
$ws.Range("H137").Value = 2469.8408
$ws.Range("I137").Value = 1398.68
$ws.Range("K137").Value = 4196.04
$ws.Range("M137").Value = -1646.04
$ws.Range("H141").Value = 1726.4762
$ws.Range("I141").Value = 1394.5294
$ws.Range("K141").Value = 4183.5882
$ws.Range("M141").Value = 996.4117999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3485.117
$ws.Range("I32").Value = 2194.4795
$ws.Range("K32").Value = 2194.4795
$ws.Range("M32").Value = -1907.4795
$ws.Range("H61").Value = 3212.611
$ws.Range("I61").Value = 2192.182
$ws.Range("K61").Value = 2192.182
$ws.Range("M61").Value = -1980.182
$ws.Range("H132").Value = 1745.0646
$ws.Range("I132").Value = 1233.8214
$ws.Range("K132").Value = 3701.4642
$ws.Range("M132").Value = -1171.4642
$ws.Range("H136").Value = 3212.611
$ws.Range("I136").Value = 2192.182
$ws.Range("K136").Value = 6576.545999999999
$ws.Range("M136").Value = -4026.545999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 169.71428
$ws.Range("I22").Value = 137.6
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 137.6
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 35.40000000000001
$ws.Range("N22").Value = -596
$ws.Range("H134").Value = 2490.45
$ws.Range("I134").Value = 783.7879
$ws.Range("K134").Value = 2351.3637
$ws.Range("M134").Value = 183.6363000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3039.8823
$ws.Range("I31").Value = 2081.9583
$ws.Range("J31").Value = 5338.9
$ws.Range("K31").Value = 2081.9583
$ws.Range("L31").Value = 5338.9
$ws.Range("M31").Value = -1786.9583
$ws.Range("N31").Value = -5928.9
$ws.Range("H34").Value = 3039.8823
$ws.Range("I34").Value = 2081.9583
$ws.Range("J34").Value = 5338.9
$ws.Range("K34").Value = 2081.9583
$ws.Range("L34").Value = 5338.9
$ws.Range("M34").Value = -1879.9583
$ws.Range("N34").Value = -5742.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1658.6316
$ws.Range("I5").Value = 1152.8
$ws.Range("J5").Value = 2220.6667
$ws.Range("K5").Value = 3458.4
$ws.Range("L5").Value = 6662.000100000001
$ws.Range("M5").Value = -3346.4
$ws.Range("N5").Value = -6886.000100000001
$ws.Range("H33").Value = 4860.2856
$ws.Range("I33").Value = 84.92856999999999
$ws.Range("K33").Value = 509.57142
$ws.Range("M33").Value = -226.57142
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H76").Value = 82899.8
$ws.Range("J76").Value = 3999.6667
$ws.Range("L76").Value = 11999.0001
$ws.Range("N76").Value = -12765.0001
$ws.Range("H79").Value = 82899.8
$ws.Range("J79").Value = 3999.6667
$ws.Range("L79").Value = 11999.0001
$ws.Range("N79").Value = -14651.0001
$ws.Range("H113").Value = 2723.7778
$ws.Range("I113").Value = 4144.6924
$ws.Range("J113").Value = 1920.6522
$ws.Range("K113").Value = 12434.0772
$ws.Range("L113").Value = 5761.9566
$ws.Range("M113").Value = -10264.0772
$ws.Range("N113").Value = -10101.9566
$ws.Range("H122").Value = 1296.8572
$ws.Range("I122").Value = 1255.8
$ws.Range("K122").Value = 11302.2
$ws.Range("M122").Value = -8852.199999999999
$ws.Range("H135").Value = 1658.6316
$ws.Range("I135").Value = 1152.8
$ws.Range("J135").Value = 2220.6667
$ws.Range("K135").Value = 10375.2
$ws.Range("L135").Value = 19986.0003
$ws.Range("M135").Value = -7840.199999999999
$ws.Range("N135").Value = -25056.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2733.6
$ws.Range("I113").Value = 1361.2
$ws.Range("J113").Value = 4106
$ws.Range("K113").Value = 1361.2
$ws.Range("L113").Value = 4106
$ws.Range("M113").Value = 808.8
$ws.Range("N113").Value = -8446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6266.4443
$ws.Range("I7").Value = 3680
$ws.Range("K7").Value = 3680
$ws.Range("M7").Value = -3568
$ws.Range("H68").Value = 3082.4285
$ws.Range("I68").Value = 3194.25
$ws.Range("J68").Value = 2933.3333
$ws.Range("K68").Value = 3194.25
$ws.Range("L68").Value = 2933.3333
$ws.Range("M68").Value = -2445.25
$ws.Range("N68").Value = -4431.3333
$ws.Range("H71").Value = 3082.4285
$ws.Range("I71").Value = 3194.25
$ws.Range("J71").Value = 2933.3333
$ws.Range("K71").Value = 15971.25
$ws.Range("L71").Value = 14666.6665
$ws.Range("M71").Value = -12227.25
$ws.Range("N71").Value = -22154.6665
$ws.Range("H122").Value = 7073.778
$ws.Range("I122").Value = 3870.5
$ws.Range("J122").Value = 9636.4
$ws.Range("K122").Value = 11611.5
$ws.Range("L122").Value = 28909.2
$ws.Range("M122").Value = -9161.5
$ws.Range("N122").Value = -33809.2
$ws.Range("H126").Value = 6266.4443
$ws.Range("I126").Value = 3680
$ws.Range("K126").Value = 11040
$ws.Range("M126").Value = -8570
$ws.Range("H132").Value = 6859.6665
$ws.Range("I132").Value = 7419.8276
$ws.Range("K132").Value = 22259.4828
$ws.Range("M132").Value = -19729.4828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5458.857
$ws.Range("I132").Value = 5850.407
$ws.Range("K132").Value = 17551.221
$ws.Range("M132").Value = -15021.221
